$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.202.36'
$ws.Range('E2').Value = '  +1.03%  '
$ws.Range('D3').Value = '2.394.81'
$ws.Range('E3').Value = '  +6.40%  '
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('D5').Value = "'329.61"
$ws.Range('E5').Value = '  +11.75%  '
$ws.Range('D6').Value = "'105.55"
$ws.Range('E6').Value = '  -5.87%  '
$ws.Range('E7').Value = '  +2.82%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = "'0.651"
$ws.Range('E9').Value = '  +7.18%  '
$ws.Range('D10').Value = "'42.01"
$ws.Range('E10').Value = '  -4.48%  '
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('E12').Value = '  -2.52%  '
$ws.Range('E13').Value = '  -1.43%  '
$ws.Range('D14').Value = "'17.17"
$ws.Range('E14').Value = '  +12.71%  '
$ws.Range('E15').Value = '  +1.76%  '
$ws.Range('D16').Value = '2.756.41'
$ws.Range('E16').Value = '  +6.44%  '
$ws.Range('D17').Value = '2.394.86'
$ws.Range('E17').Value = '  +4.98%  '
$ws.Range('D18').Value = '43.188.32'
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range('D19').Value = "'7.75"
$ws.Range('E19').Value = '  +7.84%  '
$ws.Range('E20').Value = '  +1.85%  '
$ws.Range('D21').Value = "'76.48"
$ws.Range('E21').Value = '  +1.75%  '
$ws.Range('D22').Value = "'3.73"
$ws.Range('E22').Value = '  +7.59%  '
$ws.Range('D23').Value = "'272.50"
$ws.Range('E23').Value = '  +6.64%  '
$ws.Range('E24').Value = '  -1.20%  '
$ws.Range('D25').Value = "'9.59"
$ws.Range('E25').Value = '  +7.51%  '
$ws.Range('D26').Value = "'11.81"
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').Value = "'23.01"
$ws.Range('E28').Value = '  +3.55%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'2.19"
$ws.Range('E29').Value = '  -1.91%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = "'175.34"
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').Value = "'37.60"
$ws.Range('E31').Value = '  -1.20%  '
$ws.Range('D32').Value = "'3.17"
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('D33').Value = "'0.0930"
$ws.Range('E33').Value = '  +4.61%  '
$ws.Range('D34').Value = "'5.92"
$ws.Range('E34').Value = '  +4.22%  '
$ws.Range('E35').Value = '  +4.95%  '
$ws.Range('D36').Value = "'4.94"
$ws.Range('E36').Value = '  -2.43%  '
$ws.Range('D37').Value = "'4.15"
$ws.Range('E37').Value = '  -1.74%  '
$ws.Range('E38').Value = '  -2.97%  '
$ws.Range('E39').Value = '  +3.95%  '
$ws.Range('D40').Value = "'2.82"
$ws.Range('E40').Value = '  +17.03%  '
$ws.Range('E41').Value = '  +19.42%  '
$ws.Range('E42').Value = '  +1.02%  '
$ws.Range('D43').Value = "'69.98"
$ws.Range('E43').Value = '  -2.67%  '
$ws.Range('D44').Value = "'122.00"
$ws.Range('E44').Value = '  +14.79%  '
$ws.Range('E45').Value = '  +0.15%  '
$ws.Range('D46').Value = "'12.43"
$ws.Range('E46').Value = '  -0.61%  '
$ws.Range('D47').Value = "'90.12"
$ws.Range('E47').Value = '  +47.75%  '
$ws.Range('D48').Value = "'9.37"
$ws.Range('E48').Value = '  +8.42%  '
$ws.Range('E49').Value = '  +0.50%  '
$ws.Range('E50').Value = '  +1.75%  '
$ws.Range('D51').Value = "'0.493"
$ws.Range('E51').Value = '  +12.30%  '
